$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 7
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 13

$ws.Range("E6").Select()
